$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing word2vec TransE evaluation row (row 6)
$ws.Range("B6").Value = "transE"
$ws.Range("C6").Value = "Word2Vec word embeddings"
$ws.Range("D6").Value = "pipeline default "
$ws.Range("E6").Value = 0.94776099999999996
$ws.Range("F6").Value = 4.7711439999999996
$ws.Range("G6").Value = 0.33170300000000003

# Clear the now-unused placeholder cell in F7
$ws.Range("F7").Clear()

# Move the active selection to E7
$ws.Range("E7").Select() | Out-Null
